$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Correct the "Bushaltestelle?" (column D) values for rows 4-11 from 1 to 0
# (Fehler mit SoC>100% korrigiert)
$ws.Range("D4:D11").Value = 0

# Update active selection to F12
$ws.Range("F12").Select()
